$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row above the current last row (row 30), pushing the
#    existing last row (and everything below it) down by one. Copy the
#    formatting of the row above (row 29, a "middle" table row) onto the
#    newly inserted row so the bottom-border styling stays on the new last
#    row (the row that was row 30, now row 31).
# ---------------------------------------------------------------------------
$ws.Rows("30:30").Insert(-4121)

$fmtSrc = $ws.Range("B29:J29")
$fmtDst = $ws.Range("B30:J30")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) New row 30 becomes a normal data row for LEANDRO DIAZ ARRIETA / period
#    2108 (the period previously shown on the old last row).
# ---------------------------------------------------------------------------
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1047374881"
$ws.Range("D30").Value = "LEANDRO DIAZ ARRIETA"
$ws.Range("E30").Value = "2108"
$ws.Range("F30").Value = 36341
$ws.Range("G30").Value = 908526

# ---------------------------------------------------------------------------
# 3) Row 31 (the former last row, shifted down, keeps the bottom-border
#    style) now holds the new worker.
# ---------------------------------------------------------------------------
$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1143407752"
$ws.Range("D31").Value = "LUIS CARLOS JIMENEZ CASTILLO"
$ws.Range("E31").Value = "2112"
$ws.Range("F31").Value = 36341
$ws.Range("G31").Value = 1000001

# ---------------------------------------------------------------------------
# 4) The periods listed for the original worker (rows 16-29) are now shown
#    in reverse chronological order.
# ---------------------------------------------------------------------------
$periods = @("2210","2209","2208","2207","2206","2205","2204","2203","2202","2201","2112","2111","2110","2109")
for ($i = 0; $i -lt $periods.Length; $i++) {
  $r = 16 + $i
  $ws.Range("E" + $r).Value = $periods[$i]
}

# Row 16's "Valor Mora" amount changes along with the reordering.
$ws.Range("F16").Value = 32707

# ---------------------------------------------------------------------------
# 5) Header totals: one more worker, higher total overdue amount.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 577822
$ws.Range("C13").Value = 2

# ---------------------------------------------------------------------------
# 6) Column D needs to widen to fit the new (longer) worker name.
# ---------------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 29.5
